$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.140781164169312
$ws.Range("B1").Value = 2.186646461486816
$ws.Range("C1").Value = 2.926440238952637
$ws.Range("D1").Value = 1.395778298377991
$ws.Range("E1").Value = 1.020025849342346
